$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.006795699771643
$ws.Range("E2").Value = 1.009532559647175
$ws.Range("F2").Value = 1.004972426319249
$ws.Range("J2").Value = 1.012071080893318
$ws.Range("L2").Value = 1.012405490914704
$ws.Range("M2").Value = 1.007859237321635
$ws.Range("N2").Value = 1.013508337914362
$ws.Range("C3").Value = 1.008197251359901
$ws.Range("E3").Value = 1.010735825058346
$ws.Range("F3").Value = 1.007032313929777
$ws.Range("J3").Value = 1.013101395776095
$ws.Range("L3").Value = 1.01341151132103
$ws.Range("M3").Value = 1.009718434850895
$ws.Range("N3").Value = 1.014540115962452
$ws.Range("C4").Value = 1.009102639152691
$ws.Range("E4").Value = 1.011513358573906
$ws.Range("F4").Value = 1.008363069200797
$ws.Range("J4").Value = 1.013766239260105
$ws.Range("L4").Value = 1.014060866618907
$ws.Range("M4").Value = 1.010918994503387
$ws.Range("N4").Value = 1.015205903600467
$ws.Range("C5").Value = 1.009482911626778
$ws.Range("E5").Value = 1.011839986557137
$ws.Range("F5").Value = 1.008922028095641
$ws.Range("J5").Value = 1.014045305955707
$ws.Range("L5").Value = 1.014333476974628
$ws.Range("M5").Value = 1.011423136912127
$ws.Range("N5").Value = 1.015485366602786
$ws.Range("C6").Value = 1.009546740631774
$ws.Range("E6").Value = 1.011894814532448
$ws.Range("F6").Value = 1.009015851539715
$ws.Range("J6").Value = 1.014092137257115
$ws.Range("L6").Value = 1.014379227428651
$ws.Range("M6").Value = 1.011507751549646
$ws.Range("N6").Value = 1.015532264410014
$ws.Range("C7").Value = 1.009107721746463
$ws.Range("E7").Value = 1.011517723955895
$ws.Range("F7").Value = 1.008370539941099
$ws.Range("J7").Value = 1.013769969858014
$ws.Range("L7").Value = 1.014064510731977
$ws.Range("M7").Value = 1.010925733108616
$ws.Range("N7").Value = 1.015209639496253
$ws.Range("C8").Value = 1.007269678681445
$ws.Range("E8").Value = 1.00993943221331
$ws.Range("F8").Value = 1.005669027589002
$ws.Range("J8").Value = 1.012419665540274
$ws.Range("L8").Value = 1.012745816713894
$ws.Range("M8").Value = 1.008488082817076
$ws.Range("N8").Value = 1.0138574175915
$ws.Range("C9").Value = 1.004018848159921
$ws.Range("E9").Value = 1.007149868128996
$ws.Range("F9").Value = 1.000891439981142
$ws.Range("J9").Value = 1.010025865257892
$ws.Range("L9").Value = 1.010409523095615
$ws.Range("M9").Value = 1.00417296713567
$ws.Range("N9").Value = 1.011460217838144
$ws.Range("C10").Value = 1.001843021329064
$ws.Range("E10").Value = 1.005284085024044
$ws.Range("F10").Value = 0.9976936178203026
$ws.Range("J10").Value = 1.008419891959954
$ws.Range("L10").Value = 1.008843141418668
$ws.Range("M10").Value = 1.001281905842057
$ws.Range("N10").Value = 1.009851963873915
$ws.Range("C11").Value = 1.000898696951886
$ws.Range("E11").Value = 1.004474645911795
$ws.Range("F11").Value = 0.9963056270408931
$ws.Range("J11").Value = 1.00772199414136
$ws.Range("L11").Value = 1.008162693502324
$ws.Range("M11").Value = 1.000026403695263
$ws.Range("N11").Value = 1.009153074960369
$ws.Range("C12").Value = 1.000547594827383
$ws.Range("E12").Value = 1.004173744240184
$ws.Range("F12").Value = 0.9957895443653951
$ws.Range("J12").Value = 1.007462379703838
$ws.Range("L12").Value = 1.007909607544663
$ws.Range("M12").Value = 0.9995594842477437
$ws.Range("N12").Value = 1.008893091840568
$ws.Range("C13").Value = 1.000622922897104
$ws.Range("E13").Value = 1.004238299660312
$ws.Range("F13").Value = 0.9959002698105659
$ws.Range("J13").Value = 1.007518085399739
$ws.Range("L13").Value = 1.00796391072889
$ws.Range("M13").Value = 0.9996596661921247
$ws.Range("N13").Value = 1.008948876644947
$ws.Range("C14").Value = 1.000869681699599
$ws.Range("E14").Value = 1.004449778228142
$ws.Range("F14").Value = 0.9962629782525636
$ws.Range("J14").Value = 1.007700542215689
$ws.Range("L14").Value = 1.008141780263488
$ws.Range("M14").Value = 0.9999878197400793
$ws.Range("N14").Value = 1.009131592570504
$ws.Range("C15").Value = 1.001021672823927
$ws.Range("E15").Value = 1.004580045112463
$ws.Range("F15").Value = 0.9964863851390769
$ws.Range("J15").Value = 1.007812908730181
$ws.Range("L15").Value = 1.00825132662824
$ws.Range("M15").Value = 1.000189929797714
$ws.Range("N15").Value = 1.009244118658334
$ws.Range("C16").Value = 1.001905646230184
$ws.Range("E16").Value = 1.005337771633596
$ws.Range("F16").Value = 0.9977856624107103
$ws.Range("J16").Value = 1.008466155736214
$ws.Range("L16").Value = 1.008888253624534
$ws.Range("M16").Value = 1.00136515060923
$ws.Range("N16").Value = 1.009898293350044
$ws.Range("C17").Value = 1.002459548863428
$ws.Range("E17").Value = 1.005812654693122
$ws.Range("F17").Value = 0.9985997613745494
$ws.Range("J17").Value = 1.008875244892429
$ws.Range("L17").Value = 1.00928718825863
$ws.Range("M17").Value = 1.00210134256527
$ws.Range("N17").Value = 1.010307963459785
$ws.Range("C18").Value = 1.002782421593639
$ws.Range("E18").Value = 1.006089497654568
$ws.Range("F18").Value = 0.9990742934944384
$ws.Range("J18").Value = 1.009113618959026
$ws.Range("L18").Value = 1.009519668882094
$ws.Range("M18").Value = 1.002530399804981
$ws.Range("N18").Value = 1.010546676044897
$ws.Range("C19").Value = 1.002892477678187
$ws.Range("E19").Value = 1.00618386904272
$ws.Range("F19").Value = 0.999236043444976
$ws.Range("J19").Value = 1.009194857778418
$ws.Range("L19").Value = 1.009598903173332
$ws.Range("M19").Value = 1.002676638324288
$ws.Range("N19").Value = 1.010628030232731
$ws.Range("C20").Value = 1.002400142081356
$ws.Range("E20").Value = 1.005761719632777
$ws.Range("F20").Value = 0.9985124493247256
$ws.Range("J20").Value = 1.008831378456744
$ws.Range("L20").Value = 1.009244408248378
$ws.Range("M20").Value = 1.002022392601343
$ws.Range("N20").Value = 1.010264034728729
$ws.Range("C21").Value = 1.000797026801522
$ws.Range("E21").Value = 1.004387509750593
$ws.Range("F21").Value = 0.996156184230833
$ws.Range("J21").Value = 1.007646823902335
$ws.Range("L21").Value = 1.00808941149181
$ws.Range("M21").Value = 0.9998912026063485
$ws.Range("N21").Value = 1.009077797970982
$ws.Range("C22").Value = 0.9997871235086685
$ws.Range("E22").Value = 1.003522097095155
$ws.Range("F22").Value = 0.9946716801406039
$ws.Range("J22").Value = 1.006899820848372
$ws.Range("L22").Value = 1.007361263806506
$ws.Range("M22").Value = 0.9985479300973373
$ws.Range("N22").Value = 1.008329734087001
$ws.Range("C23").Value = 1.000322681920053
$ws.Range("E23").Value = 1.003981003448065
$ws.Range("F23").Value = 0.995458938625381
$ws.Range("J23").Value = 1.007296035183665
$ws.Range("L23").Value = 1.007747456589261
$ws.Range("M23").Value = 0.9992603450241444
$ws.Range("N23").Value = 1.008726511092096
$ws.Range("C24").Value = 1.002426986095815
$ws.Range("E24").Value = 1.005784735452404
$ws.Range("F24").Value = 0.9985519028639321
$ws.Range("J24").Value = 1.008851200554655
$ws.Range("L24").Value = 1.009263739346226
$ws.Range("M24").Value = 1.002058067776039
$ws.Range("N24").Value = 1.010283884976293
$ws.Range("C25").Value = 1.004860741150578
$ws.Range("E25").Value = 1.007872076959016
$ws.Range("F25").Value = 1.002128720824314
$ws.Range("J25").Value = 1.010646469523226
$ws.Range("L25").Value = 1.011015042468879
$ws.Range("M25").Value = 1.005290970498298
$ws.Range("N25").Value = 1.012081703432719

Write-Output "Updated vm_pu values for 380 kV case"
